$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Overview")
$ws1.Range("E3").Value = "Ready for handoff"
$ws1.Range("F3").Value = "Ready for handoff"
$ws1.Range("G3").Value = "2016-09-04 02:16:34"

$ws2 = $wb.Worksheets.Item("zh-cn")
$ws2.Range("C3").Value = "Ready for handoff"
$ws2.Range("E3").Value = "mt"
$ws2.Range("H3").Value = "2016-09-04 02:16:29"

$ws3 = $wb.Worksheets.Item("de-de")
$ws3.Range("C3").Value = "Ready for handoff"
$ws3.Range("E3").Value = "mt"
$ws3.Range("H3").Value = "2016-09-04 02:16:34"
